# Loan RBI, Variable Instalments
#
# The "Repayment schedule" sheet gets a new (blank) column inserted
# before column N, pushing the old N/O/P ("Late" / "heading" / "Outstanding")
# columns one place to the right (-> O/P/Q). This makes room for an
# extra "Variable Instalments" related column in the schedule table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a brand-new column at N; everything from N onward (N,O,P) shifts
# right by one (-> O,P,Q), and the new N column comes in blank.
$ws.Columns("N:N").Insert()

# The freshly inserted column inherits formatting from its neighbour, but
# Excel still needs an explicit width for it once data columns have been
# shuffled. ColumnWidth is expressed in "characters" and gets padded by
# Excel into the internal xlsx width unit, so we dial it in so the stored
# width lands on 11 (matching column M's width).
$ws.Columns("N:N").ColumnWidth = 10.166666666666666

# Leave the cursor where it ended up after the insert.
$null = $ws.Range("R6").Select()
